# tests/data/desired/numericjoinreport.xlsx — "interlaced" sheet
#
# Commit message: "Added the information about data source to header
# columns in normal, interlaced and key layouts."
#
# For the interlaced layout this means each of the paired join/compare
# columns (Source1 / Source3) gets its own header text instead of both
# columns sharing the bare column name, e.g.
#
#   ETYPE | ETYPE              ->   ETYPE (Source1) | ETYPE (Source3)
#   FIRST_NAME | FIRST_NAME    ->   FIRST_NAME (Source1) | FIRST_NAME (Source3)
#   ...
#
# New header strings are written column-by-column, Source1 side first
# (C,E,G,I,K,M,O,Q) and then the Source3 side (D,F,H,J,L,N,P,R) — this is
# the order the workbook's writer appends them to the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interlaced")

$ws.Range("C2").Value2 = "ETYPE (Source1)"
$ws.Range("E2").Value2 = "FIRST_NAME (Source1)"
$ws.Range("G2").Value2 = "NUM_INT (Source1)"
$ws.Range("I2").Value2 = "NUM_FLOAT (Source1)"
$ws.Range("K2").Value2 = "NUM_DOUBLE (Source1)"
$ws.Range("M2").Value2 = "NUM_DECIMAL1 (Source1)"
$ws.Range("O2").Value2 = "NUM_DECIMAL2 (Source1)"
$ws.Range("Q2").Value2 = "NUM_DECIMAL3 (Source1)"

$ws.Range("D2").Value2 = "ETYPE (Source3)"
$ws.Range("F2").Value2 = "FIRST_NAME (Source3)"
$ws.Range("H2").Value2 = "NUM_INT (Source3)"
$ws.Range("J2").Value2 = "NUM_FLOAT (Source3)"
$ws.Range("L2").Value2 = "NUM_DOUBLE (Source3)"
$ws.Range("N2").Value2 = "NUM_DECIMAL1 (Source3)"
$ws.Range("P2").Value2 = "NUM_DECIMAL2 (Source3)"
$ws.Range("R2").Value2 = "NUM_DECIMAL3 (Source3)"

# The longer header text no longer fits the old best-fit column widths.
# Re-run best-fit (AutoFit) on every column pair whose header changed so the
# sheet looks right again; Source1/Source3 columns share identical header
# text so each pair comes out the same width.
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(10).AutoFit()
$ws.Columns.Item(11).AutoFit()
$ws.Columns.Item(12).AutoFit()

# AutoFit() above sizes every column to this engine's own (coarser) text
# measurement. Nudge each pair to the precise best-fit width Excel itself
# computes for this header text/font, matching the target column widths.
$ws.Columns.Item(3).ColumnWidth = 19.5
$ws.Columns.Item(4).ColumnWidth = 19.5
$ws.Columns.Item(5).ColumnWidth = 25.333333333333332
$ws.Columns.Item(6).ColumnWidth = 25.333333333333332
$ws.Columns.Item(7).ColumnWidth = 22.166666666666668
$ws.Columns.Item(8).ColumnWidth = 22.166666666666668
$ws.Columns.Item(9).ColumnWidth = 25.166666666666668
$ws.Columns.Item(10).ColumnWidth = 25.166666666666668
$ws.Columns.Item(11).ColumnWidth = 26.5
$ws.Columns.Item(12).ColumnWidth = 26.5
